$d = $word.ActiveDocument

function Split-PlainField($paraIndex, $expectedText, $word1) {
    # Replaces the paragraph's text (expected to equal $expectedText) with
    # "`" + $word1 + "~" rendered as three separate, identically-formatted
    # runs. This mirrors the target OOXML, which stores the "field-marker"
    # text as three discrete <w:r> elements rather than one merged run.
    #
    # Technique: after writing the full replacement text into the
    # paragraph (which collapses it into a single run), we temporarily drop
    # a bookmark at each internal run boundary and immediately delete it
    # again. The runtime keeps the run split that the bookmark insertion
    # produced even after the bookmark itself is removed.
    $p = $d.Paragraphs($paraIndex)
    $rng = $p.Range
    # Exclude the trailing paragraph mark from the text comparison/replace.
    $full = $d.Range($rng.Start, $rng.End - 1)
    if ($full.Text -ne $expectedText) {
        throw "Paragraph $paraIndex text was '" + $full.Text + "', expected '" + $expectedText + "'"
    }

    $full.Text = "``" + $word1 + "~"

    $boundary1 = $d.Range($full.Start + 1, $full.Start + 1)
    $d.Bookmarks.Add("zzTempSplitA", $boundary1)

    $boundary2 = $d.Range($full.Start + 1 + $word1.Length, $full.Start + 1 + $word1.Length)
    $d.Bookmarks.Add("zzTempSplitB", $boundary2)

    $d.Bookmarks("zzTempSplitA").Delete()
    $d.Bookmarks("zzTempSplitB").Delete()
}

# --- Title page (page 1) ---

# Paragraph 2: "Short Title" -> "`" / "SHORT" / (moved _GoBack bookmark) / "~"
# The _GoBack bookmark previously sat alone in paragraph 1; it is relocated
# here, between the "SHORT" and "~" runs, exactly as in the target markup.
$p2 = $d.Paragraphs(2)
$rng2 = $p2.Range
$full2 = $d.Range($rng2.Start, $rng2.End - 1)
if ($full2.Text -ne "Short Title") {
    throw "Paragraph 2 text was '" + $full2.Text + "', expected 'Short Title'"
}
$full2.Text = "``SHORT~"

$boundaryA = $d.Range($full2.Start + 1, $full2.Start + 1)
$d.Bookmarks.Add("zzTempSplitA", $boundaryA)

$goBackPos = $full2.Start + 1 + 5
$boundaryGoBack = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $boundaryGoBack)

$d.Bookmarks("zzTempSplitA").Delete()

# Paragraph 3: "Address" -> "`" / "ADDRESS" / "~"
Split-PlainField 3 "Address" "ADDRESS"

# Paragraph 6: "Section" -> "`" / "SECTIONNO" / "~"
Split-PlainField 6 "Section" "SECTIONNO"

# Paragraph 9: "09/14/2015" -> "`" / "DATE" / "~"
Split-PlainField 9 "09/14/2015" "DATE"

# Remaining "Short Title" occurrences (one per following section title block)
Split-PlainField 12 "Short Title" "SHORT"
Split-PlainField 23 "Short Title" "SHORT"
Split-PlainField 35 "Short Title" "SHORT"
Split-PlainField 47 "Short Title" "SHORT"
Split-PlainField 59 "Short Title" "SHORT"
Split-PlainField 71 "Short Title" "SHORT"
Split-PlainField 83 "Short Title" "SHORT"
Split-PlainField 95 "Short Title" "SHORT"
